$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.008791666666666666
$ws.Range("H2").Value = 0.026375
$ws.Range("I2").Value = [double]"6.529429601061531E-05"
$ws.Range("J2").Value = [double]"6.529429601061531E-05"
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.008791666666666666
$ws.Range("N2").Value = 0.026375
$ws.Range("O2").Value = [double]"6.529429601061531E-05"
$ws.Range("P2").Value = [double]"6.529429601061531E-05"
$ws.Range("Q2").Value = [double]"7.729340277777777E-05"
$ws.Range("R2").Value = 0.000695640625
$ws.Range("S2").Value = [double]"4.263345091521855E-09"
$ws.Range("T2").Value = [double]"4.263345091521855E-09"

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.008791666666666666
$ws.Range("H3").Value = 0.026375
$ws.Range("I3").Value = [double]"6.529429601061531E-05"
$ws.Range("J3").Value = [double]"6.529429601061531E-05"
$ws.Range("M3").Value = 134.6379876666666
$ws.Range("N3").Value = 403.913963
$ws.Range("O3").Value = 0.9999347057039893
$ws.Range("P3").Value = 0.9999347057039895
$ws.Range("Q3").Value = 1.183692308236111
$ws.Range("R3").Value = 10.653230774125
$ws.Range("S3").Value = [double]"6.529003266552379E-05"
$ws.Range("T3").Value = [double]"6.52900326655238E-05"

# Row 4
$ws.Range("G4").Value = 134.6379876666666
$ws.Range("H4").Value = 403.913963
$ws.Range("I4").Value = 0.9999347057039893
$ws.Range("J4").Value = 0.9999347057039895
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008791666666666666
$ws.Range("N4").Value = 0.026375
$ws.Range("O4").Value = [double]"6.529429601061531E-05"
$ws.Range("P4").Value = [double]"6.529429601061531E-05"
$ws.Range("Q4").Value = 1.183692308236111
$ws.Range("R4").Value = 10.653230774125
$ws.Range("S4").Value = [double]"6.529003266552379E-05"
$ws.Range("T4").Value = [double]"6.52900326655238E-05"

# Row 5
$ws.Range("G5").Value = 134.6379876666666
$ws.Range("H5").Value = 403.913963
$ws.Range("I5").Value = 0.9999347057039893
$ws.Range("J5").Value = 0.9999347057039895
$ws.Range("M5").Value = 134.6379876666666
$ws.Range("N5").Value = 403.913963
$ws.Range("O5").Value = 0.9999347057039893
$ws.Range("P5").Value = 0.9999347057039895
$ws.Range("Q5").Value = 18127.38772292948
$ws.Range("R5").Value = 163146.4895063654
$ws.Range("S5").Value = 0.9998694156713238
$ws.Range("T5").Value = 0.999869415671324
